$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1597
$ws1.Range("G3").Value = 79
$ws1.Range("F4").Value = 5214
$ws1.Range("F5").Value = 559
$ws1.Range("F6").Value = 10359
$ws1.Range("F8").Value = 564
$ws1.Range("F9").Value = 115
$ws1.Range("F10").Value = 116
$ws1.Range("F11").Value = 821
$ws1.Range("F12").Value = 83

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 17

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1597
$ws4.Range("G3").Value = 79
$ws4.Range("F6").Value = 5214
$ws4.Range("F7").Value = 559
$ws4.Range("F8").Value = 17
$ws4.Range("F9").Value = 10359
$ws4.Range("F11").Value = 564
$ws4.Range("F12").Value = 115
$ws4.Range("F15").Value = 116
$ws4.Range("F16").Value = 821
$ws4.Range("F18").Value = 83
